$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 750
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -216
$ws.Range("N18").Value = -1568
$ws.Range("H28").Value = 12005
$ws.Range("I28").Value = 12005
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 12005
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -11520
$ws.Range("N28").ClearContents()
$ws.Range("H62").Value = 297.5
$ws.Range("I62").Value = 297.5
$ws.Range("K62").Value = 297.5
$ws.Range("M62").Value = 326.5
$ws.Range("H65").Value = 297.5
$ws.Range("I65").Value = 297.5
$ws.Range("K65").Value = 1487.5
$ws.Range("M65").Value = 1632.5
$ws.Range("H129").Value = 1665
$ws.Range("I129").Value = 1166.5714
$ws.Range("K129").Value = 3499.7142
$ws.Range("M129").Value = 1500.2858
$ws.Range("H137").Value = 5049.2666
$ws.Range("I137").Value = 3082
$ws.Range("K137").Value = 9246
$ws.Range("M137").Value = -6696
$ws.Range("H138").Value = 2766.8333
$ws.Range("I138").Value = 1169.1111
$ws.Range("J138").Value = 4364.5557
$ws.Range("K138").Value = 3507.3333
$ws.Range("L138").Value = 13093.6671
$ws.Range("M138").Value = 1632.6667
$ws.Range("N138").Value = -23373.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4088.0857
$ws.Range("I32").Value = 3767.4707
$ws.Range("K32").Value = 3767.4707
$ws.Range("M32").Value = -3480.4707
$ws.Range("H74").Value = 1281.8572
$ws.Range("I74").Value = 1299.75
$ws.Range("J74").Value = 1174.5
$ws.Range("K74").Value = 1299.75
$ws.Range("L74").Value = 1174.5
$ws.Range("M74").Value = -425.75
$ws.Range("N74").Value = -2922.5
$ws.Range("H77").Value = 1281.8572
$ws.Range("I77").Value = 1299.75
$ws.Range("J77").Value = 1174.5
$ws.Range("K77").Value = 6498.75
$ws.Range("L77").Value = 5872.5
$ws.Range("M77").Value = -2130.75
$ws.Range("N77").Value = -14608.5
$ws.Range("H102").Value = 3188.8572
$ws.Range("I102").Value = 2470.3333
$ws.Range("K102").Value = 2470.3333
$ws.Range("M102").Value = -848.3332999999998
$ws.Range("H132").Value = 2123.3572
$ws.Range("I132").Value = 1974.5264
$ws.Range("J132").Value = 3537.25
$ws.Range("K132").Value = 5923.5792
$ws.Range("L132").Value = 10611.75
$ws.Range("M132").Value = -3393.5792
$ws.Range("N132").Value = -15671.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2849.5
$ws.Range("I86").Value = 2849
$ws.Range("K86").Value = 2849
$ws.Range("M86").Value = -1726
$ws.Range("H89").Value = 2849.5
$ws.Range("I89").Value = 2849
$ws.Range("K89").Value = 14245
$ws.Range("M89").Value = -8629
$ws.Range("H135").Value = 33084.668
$ws.Range("J135").Value = 33084.668
$ws.Range("L135").Value = 33084.668
$ws.Range("N135").Value = -43224.668
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2997.5
$ws.Range("I16").Value = 2997.5
$ws.Range("K16").Value = 2997.5
$ws.Range("M16").Value = -2710.5
$ws.Range("H31").Value = 2486
$ws.Range("I31").Value = 2748.3333
$ws.Range("J31").Value = 1699
$ws.Range("K31").Value = 2748.3333
$ws.Range("L31").Value = 1699
$ws.Range("M31").Value = -2453.3333
$ws.Range("N31").Value = -2289
$ws.Range("H34").Value = 2486
$ws.Range("I34").Value = 2748.3333
$ws.Range("J34").Value = 1699
$ws.Range("K34").Value = 2748.3333
$ws.Range("L34").Value = 1699
$ws.Range("M34").Value = -2546.3333
$ws.Range("N34").Value = -2103
$ws.Range("H58").Value = 3142
$ws.Range("I58").Value = 3165.8333
$ws.Range("K58").Value = 3165.8333
$ws.Range("M58").Value = -2962.8333
$ws.Range("H86").Value = 20370.766
$ws.Range("I86").Value = 5553.3335
$ws.Range("K86").Value = 5553.3335
$ws.Range("M86").Value = -4430.3335
$ws.Range("H89").Value = 20370.766
$ws.Range("I89").Value = 5553.3335
$ws.Range("K89").Value = 27766.6675
$ws.Range("M89").Value = -22150.6675
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -3002
$ws.Range("H113").Value = 2997.5
$ws.Range("I113").Value = 2997.5
$ws.Range("K113").Value = 2997.5
$ws.Range("M113").Value = -827.5
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030
$ws.Range("H134").Value = 5997.5
$ws.Range("I134").Value = 5997.5
$ws.Range("K134").Value = 17992.5
$ws.Range("M134").Value = -15457.5
$ws.Range("H136").Value = 3142
$ws.Range("I136").Value = 3165.8333
$ws.Range("K136").Value = 9497.499899999999
$ws.Range("M136").Value = -6947.499899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 18518518
$ws.Range("I2").Value = 18518518
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 111111108
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -111110995
$ws.Range("N2").ClearContents()
$ws.Range("H34").Value = 1967.3334
$ws.Range("J34").Value = 2501.5
$ws.Range("L34").Value = 7504.5
$ws.Range("N34").Value = -7672.5
$ws.Range("H51").Value = 475
$ws.Range("J51").Value = 475
$ws.Range("L51").Value = 1425
$ws.Range("N51").Value = -2345
$ws.Range("H68").Value = 2066.6667
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 9000
$ws.Range("M68").Value = -8189
$ws.Range("H71").Value = 2066.6667
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 27000
$ws.Range("M71").Value = -22944
$ws.Range("H81").Value = 2099.5
$ws.Range("I81").Value = 2099.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6298.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5175.5
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2099.5
$ws.Range("I84").Value = 2099.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 18895.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -13279.5
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 1665
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 1997.5
$ws.Range("K92").Value = 3000
$ws.Range("L92").Value = 5992.5
$ws.Range("M92").Value = -1752
$ws.Range("N92").Value = -8488.5
$ws.Range("H98").Value = 4761.857
$ws.Range("I98").Value = 5747.5
$ws.Range("J98").Value = 4367.6
$ws.Range("K98").Value = 17242.5
$ws.Range("L98").Value = 13102.8
$ws.Range("M98").Value = -15744.5
$ws.Range("N98").Value = -16098.8
$ws.Range("H107").Value = 1287.4445
$ws.Range("J107").Value = 198.375
$ws.Range("L107").Value = 595.125
$ws.Range("N107").Value = -4435.125
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").ClearContents()
$ws.Range("H115").Value = 2028
$ws.Range("I115").Value = 2028
$ws.Range("K115").Value = 6084
$ws.Range("M115").Value = -4909
$ws.Range("H116").Value = 50000
$ws.Range("I116").Value = 50000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 150000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -146558
$ws.Range("N116").ClearContents()
